# Fix the "Iowa bug" (and similar cases for every state that doesn't use
# the RPO file): the RPO-year value in column J of the first data row of
# each state block was stuck at 20100201; it should be 10100601.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(2,6,15,20,25,30,35,40,45,50,55,60,65,73,78,82,86,91,96,109,118,123,129,134,139,144,150,155)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 10)  # column J
    if ($cell.Value2 -eq 20100201) {
        $cell.Value = 10100601
    }
}

# Update the current selection to match the saved workbook state
# (active cell A2, with A2:P158 selected).
$ws.Range("A2:P158").Select()
